$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (2-100) so stale shared strings are dropped
$ws.Range("A2:B100").ClearContents()

$data = @(
  @{Row=2; Name='126 Racecourse Road Public Housing Tower Flemington'; Val=7},
  @{Row=3; Name='3535 Opal Meadow Heights Aged Care Community Meadow Heights'; Val=26},
  @{Row=4; Name='Al Haj Halal Meats Glenroy'; Val=11},
  @{Row=5; Name='Al-Taqwa College Truganina'; Val=8},
  @{Row=6; Name='Apartment Complex 120 Racecourse Road North Melbourne'; Val=5},
  @{Row=7; Name='Australia Post Distribution Centre Sunshine West'; Val=5},
  @{Row=8; Name='Australian Lamb Colac East'; Val=5},
  @{Row=9; Name='Baxter Foods Australia Campbellfield'; Val=5},
  @{Row=10; Name='CFMEU Melbourne Office'; Val=5},
  @{Row=11; Name='CS Square Caroline Springs'; Val=11},
  @{Row=12; Name='Cafe Roco Dandenong'; Val=6},
  @{Row=13; Name='Campbellfield Ford Complex Vaccination Clinic Campbellfield'; Val=8},
  @{Row=14; Name='Cardinia Lakes Early Learning Centre Pakenham'; Val=5},
  @{Row=15; Name='Caroline Springs Police Station'; Val=8},
  @{Row=16; Name='Chemist Warehouse Campbellfield DC'; Val=5},
  @{Row=17; Name='Chemist Warehouse Fillo Drive Somerton'; Val=12},
  @{Row=18; Name='City of Wyndham Community'; Val=6},
  @{Row=19; Name='Coles Campbellfield Plaza Campbellfield'; Val=7},
  @{Row=20; Name='Coles Coburg North Village'; Val=17},
  @{Row=21; Name='Coles Pakenham Place Shopping Centre'; Val=6},
  @{Row=22; Name='Community Kids Bayswater Early Education Centre Bayswater North'; Val=16},
  @{Row=23; Name='Construction Site 1 Warde Street Footscray'; Val=5},
  @{Row=24; Name='Construction Site Olea Apartment Caulfield North'; Val=14},
  @{Row=25; Name='Costco Wholesale Epping'; Val=28},
  @{Row=26; Name='Crusader Caravans Epping'; Val=22},
  @{Row=27; Name='Dandenong Police Station Dandenong'; Val=6},
  @{Row=28; Name='DayHab Rehabilitation Treatment Centre Ringwood East'; Val=6},
  @{Row=29; Name='Direct Freight Express Campbellfield'; Val=7},
  @{Row=30; Name='Disability Residence Life without Barriers Ashwood'; Val=5},
  @{Row=31; Name='Don Watson Coldstore Derrimutg'; Val=5},
  @{Row=32; Name='Epworth Healthcare Epworth Richmond Emergency Department'; Val=6},
  @{Row=33; Name='Ermha365 Ltd Doveton'; Val=9},
  @{Row=34; Name='Fine Food Holdings Pty Ltd Dandenong South'; Val=10},
  @{Row=35; Name='Fonterra Manufacturing Workplace Campbellfield'; Val=9},
  @{Row=36; Name='General Foods Campbellfield'; Val=9},
  @{Row=37; Name='General Foods Campbellfield'; Val=11},
  @{Row=38; Name='Gladstone Parade Early Learning & Kinder Glenroy'; Val=7},
  @{Row=39; Name='Goodstart Early Learning Altona'; Val=9},
  @{Row=40; Name='Green Leaves Early Learning Cairnlea'; Val=5},
  @{Row=41; Name='Green Leaves Early Learning Centre Highlands Craigieburn'; Val=16},
  @{Row=42; Name='Greenvale Primary School'; Val=5},
  @{Row=43; Name='Hamilton Marino 236 Jasper Road McKinnon'; Val=9},
  @{Row=44; Name='ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine'; Val=10},
  @{Row=45; Name='Ibis Kingsgate Hotel Melbourne'; Val=6},
  @{Row=46; Name='Ilim College Kiewa Campus Boys Secondary Dallas'; Val=5},
  @{Row=47; Name='Industrial Galvanizers Valmont Coatings Campbellfield'; Val=14},
  @{Row=48; Name='Inghams Enterprises Thomastown'; Val=6},
  @{Row=49; Name='Kippers Seafood Werribee'; Val=6},
  @{Row=50; Name='Kool Kidz Childcare Narre Warren'; Val=13},
  @{Row=51; Name='Level Crossing Removal Project Lilydale Construction Site John Street'; Val=6},
  @{Row=52; Name='Lineage Logistics Laverton North'; Val=6},
  @{Row=53; Name='Linfox Somerton National Distribution Centre Somerton'; Val=10},
  @{Row=54; Name='Mecca D.C Warehouse Melbourne Airport'; Val=8},
  @{Row=55; Name='Melbourne Assessment Prison West Melbourne'; Val=7},
  @{Row=56; Name='Melbourne Metropolitan Remand Centre Ravenhall'; Val=8},
  @{Row=57; Name='Melbourne West Police Station Docklands'; Val=5},
  @{Row=58; Name='Mill Park Police Station Mill Park'; Val=8},
  @{Row=59; Name='MyCentre Childcare Broadmeadows'; Val=21},
  @{Row=60; Name='National Gallery of Victoria Melbourne'; Val=8},
  @{Row=61; Name='Nido Early School Ascot Vale'; Val=12},
  @{Row=62; Name='Nido Early School Glenroy'; Val=22},
  @{Row=63; Name='Northern Health Northern Hospital Epping Emergency Department Tier 1B'; Val=55},
  @{Row=64; Name='Northern Health The Northern Hospital Epping'; Val=16},
  @{Row=65; Name='OnQ Plumbing and Excavations Craigieburn'; Val=12},
  @{Row=66; Name='Oporto Coolaroo'; Val=8},
  @{Row=67; Name='Oscar Romero Catholic Primary School Craigieburn'; Val=5},
  @{Row=68; Name='Our Lady Help of Christian''s Primary School Brunswick East'; Val=9},
  @{Row=69; Name='Pacific Meat Thomastown'; Val=5},
  @{Row=70; Name='Private Residence Daycare Allumba Way Wollert'; Val=8},
  @{Row=71; Name='Ramsay Health Care Warrigal Private Hospital'; Val=5},
  @{Row=72; Name='Ravenhall Correctional Centre Ravenhall'; Val=8},
  @{Row=73; Name='Richmond Quarter 261-271 Bridge Road Construction Site Richmond'; Val=12},
  @{Row=74; Name='Sacca''s Fruit World Broadmeadows Central Shopping Centre'; Val=6},
  @{Row=75; Name='Sharpline Stainless Steel Coburg North'; Val=6},
  @{Row=76; Name='St Margaret''s Primary School OSHC Maribyrnong'; Val=12},
  @{Row=77; Name='St Vincents Hospital Emergency Department Melbourne'; Val=9},
  @{Row=78; Name='Tek Foods Somerton'; Val=10},
  @{Row=79; Name='The Huntly-Goornong Rail Works'; Val=6},
  @{Row=80; Name='The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'; Val=18},
  @{Row=81; Name='The Royal Melbourne Hospital Ward 6SE Parkville'; Val=21},
  @{Row=82; Name='The Royal Talbot Rehabilitation Centre Kew'; Val=10},
  @{Row=83; Name='ThorwestenCabinets Pakenham'; Val=14},
  @{Row=84; Name='Truganina Early Learning Centre Truganina'; Val=6},
  @{Row=85; Name='Unilodge College Square Student Accommodation 570 Lygon Street Carlton'; Val=6},
  @{Row=86; Name='Wallaby Childcare Wollert'; Val=18},
  @{Row=87; Name='Werribee Mercy Hospital Emergency Department'; Val=13},
  @{Row=88; Name='Western Health Footscray Hospital Emergency Department'; Val=8},
  @{Row=89; Name='Western Health Sunshine Hospital Emergency Department'; Val=9},
  @{Row=90; Name='Yarra Childcare Centre Truganina'; Val=5}
)

foreach ($item in $data) {
  $ws.Cells.Item($item.Row, 1).Value = $item.Name
  $ws.Cells.Item($item.Row, 2).Value = $item.Val
}

$ws.Range("A1:B90").EntireColumn.AutoFit() | Out-Null
